# Generate Report for Handback
#
# For both language sheets ("zh-cn" and "de-de"):
#  - Status (col C) moves from "Ready for handoff" to
#    "Handed back: in sync with en-US" for every row.
#  - Two new columns get populated for every row:
#      F = Latest Target File    (same file the source .md handed off)
#      G = Latest Handback File  (the translated .xlf handed back)
#    both rendered as hyperlinks, styled like the other file-name links.
#  - Latest Handback DateTime (col H) gets a real timestamp instead of
#    the zero-date placeholder.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

function Set-HandbackColumns {
    param($SheetName, $MdTarget, $MdDisplay, $XlfTarget, $XlfDisplay, $HandbackDateTime)

    $ws = $wb.Worksheets.Item($SheetName)

    foreach ($row in 2, 3) {
        # Status text (shared across both rows / both sheets in the source data)
        $ws.Cells.Item($row, 3).Value = $newStatus

        # F = Latest Target File -> same source markdown file that was handed off
        $fCell = $ws.Cells.Item($row, 6)
        $fCell.Value = $MdDisplay
        $ws.Hyperlinks.Add($fCell, $MdTarget, "", "", $MdDisplay) | Out-Null
        $fCell.Font.Underline = $true
        $fCell.Font.Color = 15570276

        # G = Latest Handback File -> the translated .xlf handed back
        $gCell = $ws.Cells.Item($row, 7)
        $gCell.Value = $XlfDisplay
        $ws.Hyperlinks.Add($gCell, $XlfTarget, "", "", $XlfDisplay) | Out-Null
        $gCell.Font.Underline = $true
        $gCell.Font.Color = 15570276

        # H = Latest Handback DateTime -> stamp the actual handback time
        $ws.Cells.Item($row, 8).Value = $HandbackDateTime
    }
}

Set-HandbackColumns `
    "zh-cn" `
    "https://github.com/OpenLocalizationTest/oltest/blob/8cde2e04730e20388edae1a0233d94d67c86b481/e2e/8a66a406-9b05-461d-91f2-e046387f6dd5.md" `
    "8a66a406-9b05-461d-91f2-e046387f6dd5.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/99b2cd4d199300bb9944668761864321efacd374/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8a66a406-9b05-461d-91f2-e046387f6dd5.ffc532d6874f5d683ea1cc264ef62be90ed1f2e8.zh-cn.xlf" `
    "8a66a406-9b05-461d-91f2-e046387f6dd5.ffc532d6874f5d683ea1cc264ef62be90ed1f2e8.zh-cn.xlf" `
    "2016-03-19 00:48:48"

Set-HandbackColumns `
    "de-de" `
    "https://github.com/OpenLocalizationTest/oltest/blob/8cde2e04730e20388edae1a0233d94d67c86b481/e2e/8a66a406-9b05-461d-91f2-e046387f6dd5.md" `
    "8a66a406-9b05-461d-91f2-e046387f6dd5.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/784d2f4bbfd8d2ff1897288faa1099b3aee0b37c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8a66a406-9b05-461d-91f2-e046387f6dd5.ffc532d6874f5d683ea1cc264ef62be90ed1f2e8.de-de.xlf" `
    "8a66a406-9b05-461d-91f2-e046387f6dd5.ffc532d6874f5d683ea1cc264ef62be90ed1f2e8.de-de.xlf" `
    "2016-03-19 00:48:53"

Write-Host "Handback report generated."
